$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.022.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.297.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.30"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.14"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.26%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.69"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +10.88%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.658.01"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.301.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.902.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.42%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.51%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0" + ([char]0x2083).ToString() + "0905"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +12.76%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.58"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.02"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.14"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.79"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.79"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.984.56"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.45"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.526.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.74"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.16%  "
